$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.757.25"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.36%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.805.73"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.96%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.35%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.31"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.96%  "
$ws.Range("E6").Value = "  -1.80%  "
$ws.Range("E7").Value = "  +0.33%  "
$ws.Range("E8").Value = "  -0.31%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06832"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.31"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.56%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07501"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.69%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.807.30"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.76%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.768"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.40%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6237"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.71%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.051.19"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.91%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000009282"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -6.39%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "75.76"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.42%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "28.716.89"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.48%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.479"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -6.24%  "
$ws.Range("E20").Value = "  +0.36%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "211.49"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.36%  "
$ws.Range("E22").Value = "  -1.86%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.842"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.86%  "
$ws.Range("E24").Value = "  +0.43%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.24"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.51%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.876"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.42%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1271"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.24%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.44"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.51%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.432"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.64%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.06191"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.20%  "
$ws.Range("E31").Value = "  -1.53%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.785"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.21%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.761"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.79%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.735"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.07%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.065"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.95%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6438"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.94%  "
$ws.Range("E37").Value = "  -1.68%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.718"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.21%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.586"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.82%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01706"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.31%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.145.45"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.58%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8824"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.67%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.007"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.73%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.15"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.30%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.960.95"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.89%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "60.50"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.10%  "
$ws.Range("E47").Value = "  -3.92%  "
$ws.Range("E48").Value = "  +0.99%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.364"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.13%  "
$ws.Range("E50").Value = "  -0.53%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4485"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.46%  "
